$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 87

# Column A holds date-like text (e.g. "2020-08-24") stored as plain text in the
# existing rows, not as a real Excel date. Force text formatting before
# assigning the value so Excel doesn't auto-convert the string into a date
# serial number, then drop back to the default (unformatted) style so the
# new row matches the look of the rest of the table.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2020-08-25"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = 568621
$ws.Cells.Item($row, 3).Value = 628937
$ws.Cells.Item($row, 4).Value = 80878
$ws.Cells.Item($row, 5).Value = 61450
$ws.Cells.Item($row, 6).Value = 25.76
